$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Cells shift right (old B -> C, old C -> D).
# Excel's default "shift right" insert copies formatting from the column to the
# left, which is why B2 will already pick up the wrap-text style used by A2.
$ws.Columns("B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "StatQuery"

# New query text for the new column (row 2), mirroring the existing query cell.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Uterine cancer, NOS'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Make sure the new column wraps text just like column A.
$ws.Range("B2").WrapText = $true

# Match column A's width for the new column (closest achievable value).
$ws.Columns("B").ColumnWidth = 75

# Restore the selection to A2 (matches the saved selection in the workbook).
[void]$ws.Range("A2").Select()
